$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "F2"
$ws.Range("C2").Value = "Gp9"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.04339299999999999
$ws.Range("H2").Value = 0.130179
$ws.Range("I2").Value = 0.0698021577815419
$ws.Range("J2").Value = 0.0698021577815419
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.580629
$ws.Range("N2").Value = 4.741887
$ws.Range("O2").Value = 0.9548170682350041
$ws.Range("P2").Value = 0.9548170682350042
$ws.Range("Q2").Value = 0.068588234197
$ws.Range("R2").Value = 0.617294107773
$ws.Range("S2").Value = 0.06664829164944902
$ws.Range("T2").Value = 0.06664829164944903
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "F2"
$ws.Range("C3").Value = "Gp9"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.04339299999999999
$ws.Range("H3").Value = 0.130179
$ws.Range("I3").Value = 0.0698021577815419
$ws.Range("J3").Value = 0.0698021577815419
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.074797
$ws.Range("N3").Value = 0.224391
$ws.Range("O3").Value = 0.04518293176499584
$ws.Range("P3").Value = 0.04518293176499584
$ws.Range("Q3").Value = 0.003245666220999999
$ws.Range("R3").Value = 0.029210995989
$ws.Range("S3").Value = 0.003153866132092881
$ws.Range("T3").Value = 0.003153866132092881
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "F2"
$ws.Range("C4").Value = "Gp9"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.4205383333333333
$ws.Range("H4").Value = 1.261615
$ws.Range("I4").Value = 0.6764796878879081
$ws.Range("J4").Value = 0.6764796878879081
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.580629
$ws.Range("N4").Value = 4.741887
$ws.Range("O4").Value = 0.9548170682350041
$ws.Range("P4").Value = 0.9548170682350042
$ws.Range("Q4").Value = 0.6647150852783333
$ws.Range("R4").Value = 5.982435767505
$ws.Range("S4").Value = 0.6459143523096631
$ws.Range("T4").Value = 0.6459143523096631
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "F2"
$ws.Range("C5").Value = "Gp9"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.4205383333333333
$ws.Range("H5").Value = 1.261615
$ws.Range("I5").Value = 0.6764796878879081
$ws.Range("J5").Value = 0.6764796878879081
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.074797
$ws.Range("N5").Value = 0.224391
$ws.Range("O5").Value = 0.04518293176499584
$ws.Range("P5").Value = 0.04518293176499584
$ws.Range("Q5").Value = 0.03145500571833333
$ws.Range("R5").Value = 0.283095051465
$ws.Range("S5").Value = 0.03056533557824504
$ws.Range("T5").Value = 0.03056533557824504
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "F2"
$ws.Range("C6").Value = "Gp9"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.1577256666666667
$ws.Range("H6").Value = 0.473177
$ws.Range("I6").Value = 0.2537181543305499
$ws.Range("J6").Value = 0.2537181543305499
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.580629
$ws.Range("N6").Value = 4.741887
$ws.Range("O6").Value = 0.9548170682350041
$ws.Range("P6").Value = 0.9548170682350042
$ws.Range("Q6").Value = 0.2493057627776667
$ws.Range("R6").Value = 2.243751864999
$ws.Range("S6").Value = 0.242254424275892
$ws.Range("T6").Value = 0.2422544242758921
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "F2"
$ws.Range("C7").Value = "Gp9"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.1577256666666667
$ws.Range("H7").Value = 0.473177
$ws.Range("I7").Value = 0.2537181543305499
$ws.Range("J7").Value = 0.2537181543305499
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.074797
$ws.Range("N7").Value = 0.224391
$ws.Range("O7").Value = 0.04518293176499584
$ws.Range("P7").Value = 0.04518293176499584
$ws.Range("Q7").Value = 0.01179740668966667
$ws.Range("R7").Value = 0.106176660207
$ws.Range("S7").Value = 0.01146373005465792
$ws.Range("T7").Value = 0.01146373005465792

# Remove the now-unused trailing rows (old data had 9 rows, new data only has 6)
$ws.Range("A8:T10").Delete()
